# Trade #76 closed at 2026-02-16 21:36:09 - momentum UP +0.000%
#
# This script:
#  1. Updates the "Summary" sheet aggregate stats (row 2 = OVERALL, row 3 = leadlag)
#  2. Marks trade #52 (leadlag, row 42) as CLOSED on the "leadlag" sheet
#  3. Appends the newly-opened trade #76 (momentum) as a new row on the "momentum" sheet
#  4. Appends the now-closed trade #52 (leadlag) as a new row on the "All Trades" sheet
#  5. Updates the "Comparison" sheet leadlag row with refreshed stats
#
# Note: some text cells in this workbook look like numbers/percentages/dates
# (e.g. "67.3%", "2.90", "2026-02-16") but must stay literal text, matching
# the source file's inline-string cells. Assigning such a string straight to
# .Value lets Excel auto-coerce it into a real number/date, so for those
# cells specifically we first force the cell to Text format ("@") before
# writing the value. Plain words/sentences/times (e.g. "CLOSED", "momentum",
# "21:36:09") are unaffected by Excel's auto-detection and are written as-is.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
}

# ---------------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("C2").Value = 52
Set-TextValue $summary.Range("D2") "67.3%"
Set-TextValue $summary.Range("E2") "+13.6387%"
Set-TextValue $summary.Range("F2") "+0.2623%"

Set-TextValue $summary.Range("D3") "46.4%"
Set-TextValue $summary.Range("E3") "+9.4371%"
Set-TextValue $summary.Range("F3") "+0.1685%"

# ---------------------------------------------------------------------------
# 2) leadlag sheet - close trade #52 (row 42)
# ---------------------------------------------------------------------------
$leadlag = $wb.Worksheets.Item("leadlag")

$leadlag.Range("G42").Value = 68374.33375600001
$leadlag.Range("H42").Value = "CLOSED"
$leadlag.Range("I42").Value = 0.4778
$leadlag.Range("J42").Value = 4.78
$leadlag.Range("M42").Value = "time_exit_5min"
$leadlag.Range("N42").Value = 5

# ---------------------------------------------------------------------------
# 3) momentum sheet - append newly opened trade #76 (row 21)
# ---------------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")

$momentum.Range("A21").Value = 76
Set-TextValue $momentum.Range("B21") "2026-02-16"
$momentum.Range("C21").Value = "21:36:09"
$momentum.Range("D21").Value = "momentum"
$momentum.Range("E21").Value = "UP"
$momentum.Range("F21").Value = 68679.56
$momentum.Range("H21").Value = "OPEN"
$momentum.Range("I21").Value = 0
$momentum.Range("J21").Value = 0
$momentum.Range("K21").Value = 0.9
$momentum.Range("L21").Value = "Upward momentum: 0.287% over 10 samples"
$momentum.Range("N21").Value = 0

# ---------------------------------------------------------------------------
# 4) All Trades sheet - append trade #52 (leadlag, now CLOSED) as row 53
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("A53").Value = 52
Set-TextValue $allTrades.Range("B53") "2026-02-16"
$allTrades.Range("C53").Value = "21:31:09"
$allTrades.Range("D53").Value = "leadlag"
$allTrades.Range("E53").Value = "DOWN"
$allTrades.Range("F53").Value = 68702.565
$allTrades.Range("G53").Value = 68374.33375600001
$allTrades.Range("H53").Value = "CLOSED"
$allTrades.Range("I53").Value = 0.4778
$allTrades.Range("J53").Value = 4.78
$allTrades.Range("K53").Value = 0.75
$allTrades.Range("L53").Value = "Coinbase leading with -0.127% move"
$allTrades.Range("M53").Value = "time_exit_5min"
$allTrades.Range("N53").Value = 5

# ---------------------------------------------------------------------------
# 5) Comparison sheet - refreshed leadlag stats (row 2)
# ---------------------------------------------------------------------------
$comparison = $wb.Worksheets.Item("Comparison")

Set-TextValue $comparison.Range("C2") "46.4%"
Set-TextValue $comparison.Range("D2") "2.90"
Set-TextValue $comparison.Range("E2") "+0.5543%"
Set-TextValue $comparison.Range("G2") "1.67"
